# Simplify the Workblocks sheet: collapse wbInit/wbGetTransactionData/wbProcessTransaction
# (each previously with a _Type row + a _SuppressSuccessful row) down to four single
# _Type rows (Init/GetTransactionData/ProcessTransaction/CloseAllApplications), since
# wbHandleExec (and the "suppress successful" logging toggle) no longer exists.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workblocks")

# Remove the two now-obsolete rows (old rows 6 and 7) so we end up with 5 rows total.
$null = $ws.Rows.Item(6).Delete()
$null = $ws.Rows.Item(6).Delete()

# Row 2: Init workblock type
$ws.Cells.Item(2, 1).Value = "wbInitAllApplications_Type"
$ws.Cells.Item(2, 2).Value = "Init"
$ws.Cells.Item(2, 3).Value = "Name of Workblock"

# Row 3: GetTransactionData workblock type
$ws.Cells.Item(3, 1).Value = "wbGetTransactionData_Type"
$ws.Cells.Item(3, 2).Value = "GetData"
$ws.Cells.Item(3, 3).Value = "Name of Workblock"

# Row 4: ProcessTransaction workblock type
$ws.Cells.Item(4, 1).Value = "wbProcessTransaction_Type"
$ws.Cells.Item(4, 2).Value = "Process"
$ws.Cells.Item(4, 3).Value = "Name of Workblock"

# Row 5: new CloseAllApplications workblock type
$ws.Cells.Item(5, 1).Value = "wbCloseAllApplications_Type"
$ws.Cells.Item(5, 2).Value = "Close"
$ws.Cells.Item(5, 3).Value = "Name of Workblock"

# Move the selected cell / active sheet from Credentials to Workblocks, matching
# the new activeTab / tabSelected state.
$null = $ws.Activate()
$null = $ws.Range("B13").Select()
